$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Website" header column, styled like the existing header cells.
$ws.Range("C1").Value = "Website"
$ws.Range("C1").Style = $ws.Range("B1").Style

# Row 2: petes.com hyperlink
$ws.Range("C2").Value = "petes.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "http://petes.com")
$ws.Range("C2").Font.Underline = $true
$ws.Range("C2").Font.Color = 16711680

# Row 3: foodwars.com hyperlink
$ws.Range("C3").Value = "foodwars.com"
$ws.Hyperlinks.Add($ws.Range("C3"), "http://foodwars.com")
$ws.Range("C3").Font.Underline = $true
$ws.Range("C3").Font.Color = 16711680
